$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7183
$ws.Range("C3").Value = 163503
$ws.Range("C4").Value = 154494
$ws.Range("C8").Value = 64.78
